{"js": "// The document contains a title, an intro paragraph (\"These are the\n// orphan tags that were found in the documents: \"), and then one\n// paragraph per orphan tag (e.g. \"PUMP:RISK:10 \", \"ACE:SRS:110\", ...).\n// The edit removes every one of those orphan-tag paragraphs, leaving\n// just the title and the intro paragraph in front of the section\n// properties.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Orphan tag lines look like \"PUMP:RISK:10 \" / \"ACE:SRS:110\" - a tag\n// namespace, a colon, then more colon-separated segments. Walk back to\n// front so deleting doesn't disturb the indices of paragraphs we still\n// need to inspect.\nconst orphanTagPattern = /^\\s*[A-Za-z0-9]+:[A-Za-z0-9]+:[A-Za-z0-9]+\\s*$/;\n\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const paragraph = paragraphs.items[i];\n  if (orphanTagPattern.test(paragraph.text)) {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a title, an intro paragraph (\"These are the\n# orphan tags that were found in the documents: \"), and then one\n# paragraph per orphan tag (e.g. \"PUMP:RISK:10 \", \"ACE:SRS:110\", ...).\n# The edit removes every one of those orphan-tag paragraphs, leaving\n# just the title and the intro paragraph in front of the section\n# properties.\n\n$d = $word.ActiveDocument\n\n# Orphan tag lines look like \"PUMP:RISK:10 \" / \"ACE:SRS:110\" - a tag\n# namespace, a colon, then more colon-separated segments. Walk back to\n# front so deleting doesn't shift the index of paragraphs we still need\n# to inspect.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -match \"^\\s*[A-Za-z0-9]+:[A-Za-z0-9]+:[A-Za-z0-9]+\\s*$\") {\n        $p.Range.Delete()\n    }\n}\n"}
